$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Paragraph 1: drop the lead-in "A few months ago, I wrote about the "
#    so the paragraph starts directly with "Predictive Power Score (PPS)"
#    (that run keeps its hyperlink styling even after the hyperlink is
#    unwrapped below).
# ------------------------------------------------------------------
$d.Content.Find.Execute("A few months ago, I wrote about the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Merge the "Florian Wetschoreck ... ppscore" paragraph away: replace
#    its content with the following paragraph's text ("Yet, I work
#    mostly in R ... powertool ... workflow."), then drop the now
#    duplicated "Yet, I work..." paragraph plus the "So, over the
#    holiday period ... I wrote an R package!" paragraph and the blank
#    paragraph that followed it.
# ------------------------------------------------------------------
$florianIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Florian Wetschoreck*") {
        $florianIdx = $i
        break
    }
}

$florianPara = $d.Paragraphs.Item($florianIdx)
$yetPara = $d.Paragraphs.Item($florianIdx + 1)

$srcRange = $d.Range($yetPara.Range.Start, $yetPara.Range.End - 1)
$fmt = $srcRange.FormattedText

$dstRange = $d.Range($florianPara.Range.Start, $florianPara.Range.End - 1)
$dstRange.FormattedText = $fmt

# Drop the old "Yet, I work..." / "So, over the holiday period..." /
# blank paragraphs that now sit right after the rewritten paragraph.
$startPara = $d.Paragraphs.Item($florianIdx + 1)
$endPara = $d.Paragraphs.Item($florianIdx + 3)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 3. Unwrap every remaining hyperlink (PPS, correlation matrix,
#    Florian Wetschoreck, original blog, Python implementation) while
#    keeping their run formatting / text intact.
# ------------------------------------------------------------------
while ($d.Hyperlinks.Count -gt 0) {
    $d.Hyperlinks.Item(1).Delete()
}

Write-Output "done"
